# LZ BT titrations 20220802
# Append a new titration data row (row 93) to the CRMAccuracyData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$ws.Cells.Item(93, 1).Value   = 20220802
$ws.Cells.Item(93, 2).Value   = 2220.3086600000001
$ws.Cells.Item(93, 3).Value   = 2224.4699999999998
$ws.Cells.Item(93, 4).Formula = "=100*(B93-C93)/C93"
$ws.Cells.Item(93, 5).Value   = 180
$ws.Cells.Item(93, 6).Value   = "CRM OPENED 20220702 "

# Match the author's final selection/view state.
$ws.Range("C94").Select()
